# Auto-generated PowerShell script to apply F-column ("想去人数") updates
# across all 4 worksheets, per the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1442
$ws.Cells.Item(5, 6).Value = 6943
$ws.Cells.Item(6, 6).Value = 550
$ws.Cells.Item(8, 6).Value = 47
$ws.Cells.Item(9, 6).Value = 4617
$ws.Cells.Item(10, 6).Value = 6834
$ws.Cells.Item(11, 6).Value = 12
$ws.Cells.Item(13, 6).Value = 1407
$ws.Cells.Item(14, 6).Value = 823
$ws.Cells.Item(16, 6).Value = 14
$ws.Cells.Item(17, 6).Value = 35
$ws.Cells.Item(18, 6).Value = 1137
$ws.Cells.Item(20, 6).Value = 137
$ws.Cells.Item(22, 6).Value = 192
$ws.Cells.Item(24, 6).Value = 1074
$ws.Cells.Item(25, 6).Value = 540
$ws.Cells.Item(27, 6).Value = 1170
$ws.Cells.Item(28, 6).Value = 30
$ws.Cells.Item(29, 6).Value = 122
$ws.Cells.Item(32, 6).Value = 109
$ws.Cells.Item(33, 6).Value = 8
$ws.Cells.Item(34, 6).Value = 15
$ws.Cells.Item(35, 6).Value = 2
$ws.Cells.Item(38, 6).Value = 520
$ws.Cells.Item(39, 6).Value = 385
$ws.Cells.Item(41, 6).Value = 55
$ws.Cells.Item(42, 6).Value = 325
$ws.Cells.Item(44, 6).Value = 535
$ws.Cells.Item(47, 6).Value = 7
$ws.Cells.Item(48, 6).Value = 6

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 19
$ws.Cells.Item(3, 6).Value = 19
$ws.Cells.Item(4, 6).Value = 32
$ws.Cells.Item(7, 6).Value = 518
$ws.Cells.Item(11, 6).Value = 30
$ws.Cells.Item(12, 6).Value = 126
$ws.Cells.Item(13, 6).Value = 19
$ws.Cells.Item(33, 6).Value = 586
$ws.Cells.Item(36, 6).Value = 93
$ws.Cells.Item(39, 6).Value = 100
$ws.Cells.Item(40, 6).Value = 131

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(5, 6).Value = 846
$ws.Cells.Item(6, 6).Value = 617
$ws.Cells.Item(8, 6).Value = 1349
$ws.Cells.Item(9, 6).Value = 2113

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 19
$ws.Cells.Item(4, 6).Value = 1442
$ws.Cells.Item(6, 6).Value = 846
$ws.Cells.Item(8, 6).Value = 617
$ws.Cells.Item(9, 6).Value = 617
$ws.Cells.Item(10, 6).Value = 518
$ws.Cells.Item(11, 6).Value = 6943
$ws.Cells.Item(12, 6).Value = 550
$ws.Cells.Item(13, 6).Value = 47
$ws.Cells.Item(14, 6).Value = 4617
$ws.Cells.Item(15, 6).Value = 6834
$ws.Cells.Item(16, 6).Value = 12
$ws.Cells.Item(18, 6).Value = 1407
$ws.Cells.Item(20, 6).Value = 823
$ws.Cells.Item(22, 6).Value = 1349
$ws.Cells.Item(23, 6).Value = 2113
$ws.Cells.Item(25, 6).Value = 35
$ws.Cells.Item(26, 6).Value = 1137
$ws.Cells.Item(27, 6).Value = 137
$ws.Cells.Item(28, 6).Value = 192
$ws.Cells.Item(29, 6).Value = 1074
$ws.Cells.Item(31, 6).Value = 540
$ws.Cells.Item(32, 6).Value = 1170
$ws.Cells.Item(33, 6).Value = 122
$ws.Cells.Item(35, 6).Value = 109
$ws.Cells.Item(38, 6).Value = 15
$ws.Cells.Item(40, 6).Value = 520
$ws.Cells.Item(41, 6).Value = 586
$ws.Cells.Item(42, 6).Value = 385
$ws.Cells.Item(44, 6).Value = 93
$ws.Cells.Item(45, 6).Value = 325
$ws.Cells.Item(46, 6).Value = 535
$ws.Cells.Item(47, 6).Value = 100
